# China fossil fuel imports from Russia - "update data March 24"
# Refresh the Jan-Oct 2023 monthly rows (310-319) with revised trade-statistics
# figures and append the new Dec-2023 row (320).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revise existing monthly figures (rows 310-319) ---
$ws.Range("M310").Value = 7682555619
$ws.Range("C311").Value = 714230766
$ws.Range("F311").Value = 1596359301
$ws.Range("N311").Value = 5185216488
$ws.Range("F312").Value = 1352156465
$ws.Range("B313").Value = 5046632.909
$ws.Range("D313").Value = 1938654.655
$ws.Range("F313").Value = 1261028485
$ws.Range("G313").Value = 8834112.6510000005
$ws.Range("J313").Value = 3766131348
$ws.Range("M313").Value = 9704129565
$ws.Range("C314").Value = 674810327
$ws.Range("E314").Value = 311925456
$ws.Range("F314").Value = 1310571591
$ws.Range("J314").Value = 3320507486
$ws.Range("M314").Value = 10493989499
$ws.Range("N314").Value = 5443481794
$ws.Range("B315").Value = 5724999.9170000004
$ws.Range("C315").Value = 573025675
$ws.Range("E315").Value = 212038070
$ws.Range("F315").Value = 968573703
$ws.Range("J315").Value = 3159828252
$ws.Range("M315").Value = 8060767087
$ws.Range("N315").Value = 4261482288
$ws.Range("B316").Value = 5875458.4680000003
$ws.Range("C316").Value = 551225013
$ws.Range("D316").Value = 2260972.16
$ws.Range("E316").Value = 310355173
$ws.Range("F316").Value = 1075430990
$ws.Range("G316").Value = 9964067.3579999991
$ws.Range("I316").Value = 992845331
$ws.Range("J316").Value = 3505698584
$ws.Range("M316").Value = 10541721938
$ws.Range("N316").Value = 6004367341
$ws.Range("C317").Value = 494555109
$ws.Range("E317").Value = 362312257
$ws.Range("F317").Value = 1015558695
$ws.Range("J317").Value = 3153792997
$ws.Range("M317").Value = 8736032161
$ws.Range("N317").Value = 5426633373
$ws.Range("B318").Value = 4149185.7949999999
$ws.Range("C318").Value = 409639067
$ws.Range("E318").Value = 302658775
$ws.Range("F318").Value = 923702341
$ws.Range("G318").Value = 7628201.9850000003
$ws.Range("J318").Value = 2911655881
$ws.Range("M318").Value = 8534782775
$ws.Range("N318").Value = 5533272637
$ws.Range("B319").Value = 3636675.89
$ws.Range("C319").Value = 411866953
$ws.Range("D319").Value = 1923696.9650000001
$ws.Range("E319").Value = 313966609
$ws.Range("F319").Value = 958762597
$ws.Range("G319").Value = 7306332.5580000002
$ws.Range("I319").Value = 609571153
$ws.Range("J319").Value = 4187667201
$ws.Range("M319").Value = 8999568903
$ws.Range("N319").Value = 5545020618

# --- Append the new row 320 (2023-12-31, serial 45291) ---
# Clone formatting/styles from row 319 first, then overwrite with the new values.
$ws.Range("A319:N319").Copy() | Out-Null
$ws.Range("A320:N320").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("A320").Value = 45291
$ws.Range("B320").Value = 3389292.7050000001
$ws.Range("C320").Value = 373722796
$ws.Range("D320").Value = 2425503.2000000002
$ws.Range("E320").Value = 423520291
$ws.Range("F320").Value = 1042735863
$ws.Range("G320").Value = 7592889.1310000001
$ws.Range("H320").Value = 1607432321
$ws.Range("I320").Value = 519221118
$ws.Range("J320").Value = 5851319211
$ws.Range("K320").Value = 0
$ws.Range("L320").Value = 549402142
$ws.Range("M320").Value = 9561014006
$ws.Range("N320").Value = 5525313310

# --- Leave the active selection on F315 (matches the saved view state) ---
$ws.Range("F315").Select() | Out-Null
